# Regenerate the G ("K") column values for save_data sheet.
# This mirrors: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" — the G column (header "K") is recomputed and
# rewritten for every data row (rows 2-75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 1
    7 = 2
    8 = 1
    9 = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 2
    21 = 3
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 1
    30 = 1
    31 = 0
    32 = 1
    33 = 1
    34 = 2
    35 = 0
    36 = 1
    37 = 0
    38 = 4
    39 = 1
    40 = 0
    41 = 1
    42 = 2
    43 = 3
    44 = 2
    45 = 3
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    50 = 0
    51 = 2
    52 = 1
    53 = 2
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 1
    63 = 2
    64 = 1
    65 = 1
    66 = 0
    67 = 0
    68 = 1
    69 = 1
    70 = 2
    71 = 0
    72 = 0
    73 = 1
    74 = 1
    75 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
